# test-20.xlsx: "3Y size for all metrics"
# Fills in the previously-empty train/test metric columns on the
# "mse (3Y size)" and "r2 (3Y size)" sheets, and updates the active
# sheet / selection state to match the post-edit UI snapshot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "mse (3Y size)" sheet: fill D/E (Random Forest train/test) and
#    I/J (another model's train/test) columns for rows 6-55.
#    (M/N/O "shrinkage" columns already had data and are untouched.)
# ---------------------------------------------------------------
$wsMse3y = $wb.Worksheets.Item("mse (3Y size)")

$mse3yData = @{
    6 = @(3216.9870924425632, 25351.210120072548, 3045.670391142683, 27046.08139201935)
    7 = @(3154.973506166868, 26121.8972496977, 3080.5811506045939, 22439.472626964929)
    8 = @(3100.3716781741232, 28697.06710278113, 2901.2450577085851, 23917.951183192261)
    9 = @(3240.852323458284, 25242.038786577989, 3074.1009652660209, 22844.802552962508)
    10 = @(3307.9025003627571, 21073.26174183796, 3273.487678415961, 20245.27559866989)
    11 = @(3291.8693716142679, 24302.541663603381, 3108.0674021765408, 23408.078405199511)
    12 = @(3130.1654001511488, 25536.385920072549, 3296.526512817411, 17706.922864691649)
    13 = @(3334.504792442563, 24087.751516686811, 3229.7577689540499, 19439.272640870611)
    14 = @(3355.074118863361, 21710.85839274486, 3077.4802026602169, 24794.033934703752)
    15 = @(3221.8093696493352, 24777.486592019341, 3102.4045054111239, 22629.07936203144)
    16 = @(3396.0688364570742, 19338.57599455864, 3022.8268705562268, 21568.508714873031)
    17 = @(3154.0567938331319, 24617.994903264811, 3041.4236683494551, 23575.398734461909)
    18 = @(3245.1017642382112, 24208.241315235791, 3130.967229353083, 21726.940821765409)
    19 = @(3397.996879474003, 20916.038330229741, 3086.0629383010869, 21996.453624304711)
    20 = @(3400.3085809552608, 22067.07292660217, 3094.6658655683191, 23168.47734171704)
    21 = @(3179.3582010580408, 27018.915974365169, 3134.8409441354288, 20793.223449334939)
    22 = @(3306.9044490931069, 22842.10722273276, 3044.2625793530819, 24005.32559854897)
    23 = @(3115.8719370012091, 27090.224272430471, 3092.7339331318012, 21873.356974727929)
    24 = @(3271.5938137847638, 22284.024758041109, 3116.78410828295, 23300.082884159609)
    25 = @(3483.1729377871829, 21555.31059673518, 3252.620860278113, 18971.207692623939)
    26 = @(3187.9207625755739, 25566.761665659, 3056.872164691657, 21584.457699516319)
    27 = @(3480.8487447702541, 20618.276820918982, 3197.917487575573, 19761.29710120918)
    28 = @(3441.197816021765, 22556.80917593712, 3068.2289293833128, 22698.072636759371)
    29 = @(3259.6399957980648, 25429.238391414739, 3087.8486777206772, 21509.049047037479)
    30 = @(3372.943576904474, 22076.55189383313, 3081.8367878174122, 21031.092299274489)
    31 = @(3304.854324818622, 25691.617358162031, 3145.1498453446179, 18799.522951995161)
    32 = @(3154.557795918985, 24577.280992382101, 3047.8011408706152, 23166.169537605801)
    33 = @(3367.8476159915349, 24980.123166989109, 3135.6538781438931, 20758.261880532042)
    34 = @(3260.4142163240631, 25114.813624667469, 3077.1735257255132, 21649.48267230955)
    35 = @(3224.8080128174129, 23644.69216045949, 3075.149316263602, 24549.111314631191)
    36 = @(3309.5213876662629, 25062.871717533249, 3194.785590084643, 21330.725482708582)
    37 = @(3332.2874597339778, 23240.710217775089, 3124.9209786577981, 19845.542686094312)
    38 = @(3381.7529509673532, 21502.959778234581, 3004.0029480350659, 23979.099173397819)
    39 = @(3268.376077388149, 25027.598042321639, 3095.2343067714619, 20299.23827146312)
    40 = @(3311.4699999395411, 23277.425789117289, 3168.5234801390561, 20753.645466263599)
    41 = @(3369.1776064691649, 21483.52724244256, 3069.061685792019, 24413.13493651753)
    42 = @(3351.5865921704958, 21950.386020435311, 3054.0664550785968, 21997.308454292619)
    43 = @(3418.411057889964, 19416.532257436509, 3111.4974358524778, 22067.477231801691)
    44 = @(3358.2446576481261, 22929.82514498186, 3164.7610568923819, 21556.806998911728)
    45 = @(3415.787828506649, 22472.83811305925, 3226.5033131197092, 18173.554329987899)
    46 = @(3355.8406158403868, 23758.70939951632, 3128.3264233373638, 21767.688337968561)
    47 = @(3397.3991872732759, 21938.229162152358, 3014.7862496977018, 23640.734281015721)
    48 = @(3357.0984856711011, 22992.53917037485, 3106.2381782043522, 20609.787174002409)
    49 = @(3328.5219995465532, 26237.459122007251, 3232.6292698609432, 20846.896878476411)
    50 = @(3341.2638752720682, 22715.359543530831, 3075.7291628174121, 23374.509329866982)
    51 = @(3311.3510293228542, 22470.66266311971, 3110.3434353385728, 21763.49903083434)
    52 = @(3386.4982378174118, 20027.732994316801, 3121.5836908403862, 20026.41332720677)
    53 = @(3347.7943035368799, 22530.89984691656, 3098.2479312877872, 21452.670103869401)
    54 = @(3172.8161358827078, 24861.175624909309, 3192.224364812575, 19293.88813688029)
    55 = @(3148.306387273276, 25888.576547037479, 3203.9019245163231, 20743.61531438935)
}

foreach ($row in $mse3yData.Keys) {
    $vals = $mse3yData[$row]
    $wsMse3y.Cells.Item($row, 4).Value = $vals[0]   # D
    $wsMse3y.Cells.Item($row, 5).Value = $vals[1]   # E
    $wsMse3y.Cells.Item($row, 9).Value = $vals[2]   # I
    $wsMse3y.Cells.Item($row, 10).Value = $vals[3]  # J
}

# ---------------------------------------------------------------
# 2) "r2 (3Y size)" sheet: fill I/J columns for rows 6-55.
#    (D/E already had data and are untouched.)
# ---------------------------------------------------------------
$wsR23y = $wb.Worksheets.Item("r2 (3Y size)")

$r23yData = @{
    6 = @(0.90412555723933363, 0.42750268680290798)
    7 = @(0.90939255972905075, 0.28052477256734548)
    8 = @(0.91055416855652394, 0.27854230234701188)
    9 = @(0.90881597568581884, 0.3497551292566905)
    10 = @(0.90911387652511921, 0.35376807502779112)
    11 = @(0.90539514770878049, 0.35355065251754408)
    12 = @(0.90901591236516344, 0.3300243664886775)
    13 = @(0.9079087221562393, 0.28058117093124868)
    14 = @(0.90785397301602544, 0.39191919273441711)
    15 = @(0.90737715026868282, 0.33280952346996717)
    16 = @(0.90823028905856296, 0.34196350840261958)
    17 = @(0.90534639015244911, 0.36590502410462061)
    18 = @(0.91009047458693115, 0.25799876996976823)
    19 = @(0.90929492097599995, 0.32661407871416231)
    20 = @(0.90792311678307247, 0.3178600381297112)
    21 = @(0.90619680558169891, 0.32132652363235858)
    22 = @(0.91123213881750642, 0.37199031296290869)
    23 = @(0.91187704910795719, 0.31878276838904862)
    24 = @(0.90925812610268597, 0.29646399126604911)
    25 = @(0.906255970468544, 0.29727363757154218)
    26 = @(0.9079690878834441, 0.36849869114487532)
    27 = @(0.90866727410770942, 0.27775363877604847)
    28 = @(0.90930927887235624, 0.343630577507264)
    29 = @(0.90547266585714004, 0.34908375181311868)
    30 = @(0.9078507177160523, 0.31702361724545369)
    31 = @(0.91054646095650116, 0.30692973675814178)
    32 = @(0.90904426610177547, 0.3268303380790798)
    33 = @(0.90325833808149991, 0.37223291292966559)
    34 = @(0.90996281054771588, 0.30234599334582413)
    35 = @(0.91197101567459815, 0.29175726336685348)
    36 = @(0.91160163959341389, 0.32560207606498809)
    37 = @(0.91032619961624062, 0.25557868482060792)
    38 = @(0.90649995234227743, 0.38142632022232581)
    39 = @(0.9079017055826375, 0.37956159532344391)
    40 = @(0.90975089215012217, 0.31349251974106551)
    41 = @(0.91005514605582094, 0.34460329470773082)
    42 = @(0.90996706994976218, 0.32837263881826578)
    43 = @(0.90644311154764323, 0.37410656586051089)
    44 = @(0.91034518172196011, 0.33680249949420338)
    45 = @(0.90602196001410384, 0.33026211789630988)
    46 = @(0.90786779116412131, 0.35652337109062537)
    47 = @(0.90532477036474546, 0.40728237695587011)
    48 = @(0.90723085000838521, 0.33864142300859512)
    49 = @(0.90451943460731488, 0.35211356528402188)
    50 = @(0.90629700275671321, 0.40043933595991321)
    51 = @(0.90859010188473355, 0.31176122208155832)
    52 = @(0.90547605896457839, 0.41057111654881329)
    53 = @(0.91076257084960499, 0.34449376037067903)
    54 = @(0.90777686714993711, 0.30161677552227301)
    55 = @(0.90677725150959487, 0.38646012884856662)
}

foreach ($row in $r23yData.Keys) {
    $vals = $r23yData[$row]
    $wsR23y.Cells.Item($row, 9).Value = $vals[0]    # I
    $wsR23y.Cells.Item($row, 10).Value = $vals[1]   # J
}

# ---------------------------------------------------------------
# 3) Update view state: selection / active sheet / active tab.
#    "mse" keeps its own selection (C4:O58) untouched.
#    "mse (3Y size)" loses tabSelected, gets a new selection (R55).
#    "r2 (3Y size)" becomes the active sheet with selection Q32.
# ---------------------------------------------------------------
$wsMse3y.Activate()
$wsMse3y.Range("R55").Select()

$wsR23y.Activate()
$wsR23y.Range("Q32").Select()
